$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 2.75
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 3.5
$ws.Range("W3").Value = 7
$ws.Range("X3").Value = 12
$ws.Range("AG3").Value = 7
$ws.Range("AJ3").Value = 29
$ws.Range("AN3").Value = 4.5
$ws.Range("AO3").Value = 17
$ws.Range("AR3").Value = 101
$ws.Range("AX3").Value = 17
$ws.Range("AZ3").Value = 51

# Row 4
$ws.Range("O4").Value = 1.57
$ws.Range("P4").Value = 2.25
$ws.Range("S4").Value = 1.62
$ws.Range("T4").Value = 2.2
$ws.Range("U4").Value = 2.5
$ws.Range("V4").Value = 1.5
$ws.Range("Y4").Value = 10
$ws.Range("AC4").Value = 5.5
$ws.Range("AE4").Value = 23
$ws.Range("AF4").Value = 101
$ws.Range("AG4").Value = 9
$ws.Range("AH4").Value = 23
$ws.Range("AI4").Value = 19
$ws.Range("AK4").Value = 51
$ws.Range("AT4").Value = 2.2
$ws.Range("AW4").Value = 6.5
$ws.Range("BA4").Value = 201

# Row 6
$ws.Range("Q6").Value = 1.93
$ws.Range("R6").Value = 1.93

# Row 7
$ws.Range("H7").Value = 3
$ws.Range("I7").Value = 3.3
$ws.Range("J7").Value = 2.77
$ws.Range("L7").Value = 3.85
$ws.Range("O7").Value = 1.38
$ws.Range("P7").Value = 2.6
$ws.Range("Q7").Value = 2.1
$ws.Range("R7").Value = 1.57
$ws.Range("U7").Value = 1.82
$ws.Range("V7").Value = 1.78
$ws.Range("W7").Value = 6.7
$ws.Range("X7").Value = 10.25
$ws.Range("AA7").Value = 19.5
$ws.Range("AC7").Value = 7.5
$ws.Range("AE7").Value = 15
$ws.Range("AF7").Value = 80
$ws.Range("AG7").Value = 8.5
$ws.Range("AH7").Value = 16.5
$ws.Range("AI7").Value = 11.75
$ws.Range("AJ7").Value = 45
$ws.Range("AK7").Value = 35
$ws.Range("AL7").Value = 45
$ws.Range("AM7").Value = 700
$ws.Range("AO7").Value = 11.25
$ws.Range("AP7").Value = 19.5
$ws.Range("AQ7").Value = 45
$ws.Range("AR7").Value = 75
$ws.Range("AT7").Value = 2.42
$ws.Range("AU7").Value = 6.9
$ws.Range("AW7").Value = 5.1
$ws.Range("AX7").Value = 18.5
$ws.Range("AY7").Value = 26
$ws.Range("AZ7").Value = 100
$ws.Range("BA7").Value = 150
$ws.Range("BB7").Value = 350

# Row 8
$ws.Range("I8").Value = 2.8

# Row 13
$ws.Range("J13").Value = 2.32
$ws.Range("K13").Value = 2.02
$ws.Range("W13").Value = 5.5
$ws.Range("X13").Value = 7.1
$ws.Range("AA13").Value = 16
$ws.Range("AC13").Value = 7.7
$ws.Range("AF13").Value = 100
$ws.Range("AG13").Value = 11.25
$ws.Range("AH13").Value = 27
$ws.Range("AL13").Value = 60
$ws.Range("AN13").Value = 3.4
$ws.Range("AO13").Value = 8.75
$ws.Range("AP13").Value = 20
$ws.Range("AQ13").Value = 32
$ws.Range("BA13").Value = 200
